# adjust property of scene
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CamOffestPos / CamOffestRot values for row 2 (Demo1) and row 6 (SelectScene)
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Update the current sheet view: scroll so column E is the top-left visible column,
# and move the active selection to K7
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("K7").Select()
